$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 128 (pushing the existing 128..147 down to 130..149)
$ws.Rows("128:129").Insert()

# Fill in row 128 with the new weekly price record
$ws.Range("A128").Value = 1
$ws.Range("B128").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C128").Value = "Arica y Parinacota"
$ws.Range("D128").Value = 45135
$ws.Range("E128").Value = 15
$ws.Range("F128").Value = 100112021
$ws.Range("G128").Value = "Ají"
$ws.Range("H128").Value = "Inferno"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 160
$ws.Range("K128").Value = 11000
$ws.Range("L128").Value = 12000
$ws.Range("M128").Value = 11500
$ws.Range("N128").Value = "$/caja 15 kilos"
$ws.Range("O128").Value = "Región de Arica y Parinacota"
$ws.Range("P128").Value = 767
$ws.Range("Q128").Value = 15
$ws.Range("R128").Value = "Hortaliza"

# Fill in row 129 with the new weekly price record
$ws.Range("A129").Value = 1
$ws.Range("B129").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C129").Value = "Arica y Parinacota"
$ws.Range("D129").Value = 45135
$ws.Range("E129").Value = 15
$ws.Range("F129").Value = 100112021
$ws.Range("G129").Value = "Ají"
$ws.Range("H129").Value = "Inferno"
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 200
$ws.Range("K129").Value = 9000
$ws.Range("L129").Value = 10000
$ws.Range("M129").Value = 9500
$ws.Range("N129").Value = "$/caja 15 kilos"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 633
$ws.Range("Q129").Value = 15
$ws.Range("R129").Value = "Hortaliza"
